# "article 74 is live" -- shift the rotating blog-post carousel on row 7
# forward by one slot: the post that was showing as #71 now shows #72,
# the one showing #72 now shows #73, and the one showing #73 now shows
# the brand new #74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$I7 = $ws.Range("I7")
$E7 = $ws.Range("E7")
$C7 = $ws.Range("C7")

$I7.Value = ($I7.Text -replace "ser:\s*71", "ser: 72")
$E7.Value = ($E7.Text -replace "ser:\s*72", "ser: 73")
$C7.Value = ($C7.Text -replace "ser:\s*73", "ser: 74")
